$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 225.72
$ws.Range("I53").Value = 194.81818
$ws.Range("J53").Value = 250
$ws.Range("K53").Value = 194.81818
$ws.Range("L53").Value = 250
$ws.Range("M53").Value = 442.18182
$ws.Range("N53").Value = -1524

$ws.Range("H75").Value = 32063.2
$ws.Range("J75").Value = 32063.2
$ws.Range("L75").Value = 32063.2
$ws.Range("N75").Value = -33935.2

$ws.Range("H78").Value = 32063.2
$ws.Range("J78").Value = 32063.2
$ws.Range("L78").Value = 96189.60000000001
$ws.Range("N78").Value = -105549.6

$ws.Range("H88").Value = 6264.75
$ws.Range("J88").Value = 8660.643
$ws.Range("L88").Value = 8660.643
$ws.Range("N88").Value = -9472.643

$ws.Range("H91").Value = 6264.75
$ws.Range("J91").Value = 8660.643
$ws.Range("L91").Value = 8660.643
$ws.Range("N91").Value = -11468.643

$ws.Range("H103").Value = 913.4286
$ws.Range("I103").Value = 982.3333
$ws.Range("K103").Value = 2946.9999
$ws.Range("M103").Value = -2360.9999

$ws.Range("H138").Value = 1929.2836
$ws.Range("I138").Value = 1094.225
$ws.Range("J138").Value = 3166.4075
$ws.Range("K138").Value = 3282.675
$ws.Range("L138").Value = 9499.2225
$ws.Range("M138").Value = 1857.325
$ws.Range("N138").Value = -19779.2225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 930.0278
$ws.Range("I2").Value = 738.3103599999999
$ws.Range("K2").Value = 738.3103599999999
$ws.Range("M2").Value = -625.3103599999999

$ws.Range("H74").Value = 4882.3794
$ws.Range("I74").Value = 990.2381
$ws.Range("J74").Value = 15099.25
$ws.Range("K74").Value = 990.2381
$ws.Range("L74").Value = 15099.25
$ws.Range("M74").Value = -116.2381
$ws.Range("N74").Value = -16847.25

$ws.Range("H77").Value = 4882.3794
$ws.Range("I77").Value = 990.2381
$ws.Range("J77").Value = 15099.25
$ws.Range("K77").Value = 4951.190500000001
$ws.Range("L77").Value = 75496.25
$ws.Range("M77").Value = -583.1905000000006
$ws.Range("N77").Value = -84232.25

$ws.Range("H116").Value = 930.0278
$ws.Range("I116").Value = 738.3103599999999
$ws.Range("K116").Value = 738.3103599999999
$ws.Range("M116").Value = 1555.68964

$ws.Range("H122").Value = 1361.1333
$ws.Range("I122").Value = 959.1
$ws.Range("J122").Value = 2165.2
$ws.Range("K122").Value = 2877.3
$ws.Range("L122").Value = 6495.599999999999
$ws.Range("M122").Value = -427.3000000000002
$ws.Range("N122").Value = -11395.6

$ws.Range("H124").Value = 24216.666
$ws.Range("J124").Value = 24216.666
$ws.Range("L124").Value = 24216.666
$ws.Range("N124").Value = -34036.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 930.0278
$ws.Range("I3").Value = 738.3103599999999
$ws.Range("K3").Value = 738.3103599999999
$ws.Range("M3").Value = -624.3103599999999

$ws.Range("H94").Value = 1147.5938
$ws.Range("I94").Value = 844.1429000000001
$ws.Range("J94").Value = 1383.6111
$ws.Range("K94").Value = 844.1429000000001
$ws.Range("L94").Value = 1383.6111
$ws.Range("M94").Value = -393.1429000000001
$ws.Range("N94").Value = -2285.6111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3931.9375
$ws.Range("I16").Value = 2050.9167
$ws.Range("J16").Value = 9575
$ws.Range("K16").Value = 2050.9167
$ws.Range("L16").Value = 9575
$ws.Range("M16").Value = -1763.9167
$ws.Range("N16").Value = -10149

$ws.Range("H31").Value = 12063578
$ws.Range("I31").Value = 21277746
$ws.Range("J31").Value = 33969.223
$ws.Range("K31").Value = 21277746
$ws.Range("L31").Value = 33969.223
$ws.Range("M31").Value = -21277451
$ws.Range("N31").Value = -34559.223

$ws.Range("H34").Value = 12063578
$ws.Range("I34").Value = 21277746
$ws.Range("J34").Value = 33969.223
$ws.Range("K34").Value = 21277746
$ws.Range("L34").Value = 33969.223
$ws.Range("M34").Value = -21277544
$ws.Range("N34").Value = -34373.223

$ws.Range("H113").Value = 3931.9375
$ws.Range("I113").Value = 2050.9167
$ws.Range("J113").Value = 9575
$ws.Range("K113").Value = 2050.9167
$ws.Range("L113").Value = 9575
$ws.Range("M113").Value = 119.0832999999998
$ws.Range("N113").Value = -13915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 784.63336
$ws.Range("I5").Value = 534.5
$ws.Range("J5").Value = 1159.8334
$ws.Range("K5").Value = 1603.5
$ws.Range("L5").Value = 3479.5002
$ws.Range("M5").Value = -1491.5
$ws.Range("N5").Value = -3703.5002

$ws.Range("H68").Value = 1227.0625
$ws.Range("I68").Value = 586
$ws.Range("J68").Value = 2051.2856
$ws.Range("K68").Value = 1758
$ws.Range("L68").Value = 6153.8568
$ws.Range("M68").Value = -947
$ws.Range("N68").Value = -7775.8568

$ws.Range("H71").Value = 1227.0625
$ws.Range("I71").Value = 586
$ws.Range("J71").Value = 2051.2856
$ws.Range("K71").Value = 5274
$ws.Range("L71").Value = 18461.5704
$ws.Range("M71").Value = -1218
$ws.Range("N71").Value = -26573.5704

$ws.Range("H131").Value = 875.5
$ws.Range("I131").Value = 239.85715
$ws.Range("J131").Value = 1019.0323
$ws.Range("K131").Value = 719.5714499999999
$ws.Range("L131").Value = 3057.0969
$ws.Range("M131").Value = 4320.428550000001
$ws.Range("N131").Value = -13137.0969

$ws.Range("H135").Value = 784.63336
$ws.Range("I135").Value = 534.5
$ws.Range("J135").Value = 1159.8334
$ws.Range("K135").Value = 4810.5
$ws.Range("L135").Value = 10438.5006
$ws.Range("M135").Value = -2275.5
$ws.Range("N135").Value = -15508.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3836.875
$ws.Range("I70").Value = 3813.5715
$ws.Range("K70").Value = 3813.5715
$ws.Range("M70").Value = -3543.5715

$ws.Range("H73").Value = 3836.875
$ws.Range("I73").Value = 3813.5715
$ws.Range("K73").Value = 3813.5715
$ws.Range("M73").Value = -2877.5715

$ws.Range("H97").Value = 1114.7142
$ws.Range("I97").Value = 926
$ws.Range("J97").Value = 1190.2
$ws.Range("K97").Value = 926
$ws.Range("L97").Value = 1190.2
$ws.Range("M97").Value = -430
$ws.Range("N97").Value = -2182.2

$ws.Range("H126").Value = 1809.4546
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1925.5
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 5776.5
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -10716.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 274.08694
$ws.Range("I16").Value = 270.3
$ws.Range("K16").Value = 270.3
$ws.Range("M16").Value = -100.3

$ws.Range("H22").Value = 5775.263
$ws.Range("I22").Value = 229
$ws.Range("J22").Value = 8335.076999999999
$ws.Range("K22").Value = 229
$ws.Range("L22").Value = 8335.076999999999
$ws.Range("M22").Value = 66
$ws.Range("N22").Value = -8925.076999999999

$ws.Range("H27").Value = 5775.263
$ws.Range("I27").Value = 229
$ws.Range("J27").Value = 8335.076999999999
$ws.Range("K27").Value = 229
$ws.Range("L27").Value = 8335.076999999999
$ws.Range("M27").Value = -122
$ws.Range("N27").Value = -8549.076999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14286729
$ws.Range("I122").Value = 18182746
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 54548238
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -54545788
$ws.Range("N122").Value = -8900.0002
